$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.660.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.80%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.887.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'176.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +7.62%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.671"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.754"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.90%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +5.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'54.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.52%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.49%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.503.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.885.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.80%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'21.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.41%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'13.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.53%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -3.82%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'71.447.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'441.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.05%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -4.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'94.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.85%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -3.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'13.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.79%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.93%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -5.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'8.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +14.50%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'35.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'13.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'InjectiveProtocol"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'48.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Hedera"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.127"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +11.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'69.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.87%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'635.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.69%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.438"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'ThetaToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.28%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Dai"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.17%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.14%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Fetch.AI"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'2.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.81%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'dogwifhat"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +18.90%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -3.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -14.48%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.913.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.44%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.000280"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.78%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.71%  "
$ws.Range("E51").Style = "Normal"
